$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings are preserved exactly
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '69.266.23'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').Value = '3.735.39'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '615.75'
$ws.Range('E5').Value = '  +8.30%  '
$ws.Range('D6').Value = '192.63'
$ws.Range('E6').Value = '  +12.62%  '
$ws.Range('D7').Value = '0.642'
$ws.Range('E7').Value = '  +3.65%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.732'
$ws.Range('E9').Value = '  +4.63%  '
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('D11').Value = '60.65'
$ws.Range('E11').Value = '  +16.70%  '
$ws.Range('D12').Value = '0.0000293'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '10.61'
$ws.Range('E13').Value = '  +1.53%  '
$ws.Range('D14').Value = '4.327.77'
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').Value = '3.731.64'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('D17').Value = '19.62'
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').Value = '13.06'
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('D20').Value = '69.113.97'
$ws.Range('E20').Value = '  +2.19%  '
$ws.Range('D21').Value = '414.54'
$ws.Range('E21').Value = '  +2.96%  '
$ws.Range('D22').Value = '4.60'
$ws.Range('E22').Value = '  +3.82%  '
$ws.Range('D23').Value = '90.58'
$ws.Range('E23').Value = '  +3.77%  '
$ws.Range('D24').Value = '3.10'
$ws.Range('E24').Value = '  +2.57%  '
$ws.Range('D25').Value = '11.42'
$ws.Range('E25').Value = '  +8.48%  '
$ws.Range('D26').Value = '13.07'
$ws.Range('E26').Value = '  +3.86%  '
$ws.Range('D27').Value = '3.84'
$ws.Range('E27').Value = '  +2.34%  '
$ws.Range('D28').Value = '6.05'
$ws.Range('E28').Value = '  +1.53%  '
$ws.Range('D29').Value = '9.88'
$ws.Range('E29').Value = '  +5.18%  '
$ws.Range('D30').Value = '33.15'
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('D31').Value = '7.91'
$ws.Range('E31').Value = '  +5.64%  '
$ws.Range('D32').Value = '12.86'
$ws.Range('E32').Value = '  +3.31%  '
$ws.Range('D33').Value = '648.80'
$ws.Range('E33').Value = '  +8.91%  '
$ws.Range('E34').Value = '  +6.81%  '
$ws.Range('D35').Value = '46.41'
$ws.Range('E35').Value = '  +8.74%  '
$ws.Range('D36').Value = '67.09'
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('D37').Value = '0.0₃0843'
$ws.Range('E37').Value = '  -4.42%  '
$ws.Range('D38').Value = '0.419'
$ws.Range('E38').Value = '  +6.46%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('E41').Value = '  +5.30%  '
$ws.Range('D42').Value = '3.11'
$ws.Range('E42').Value = '  +3.98%  '
$ws.Range('D43').Value = '0.0452'
$ws.Range('E43').Value = '  +4.12%  '
$ws.Range('E44').Value = '  +4.27%  '
$ws.Range('D45').Value = '2.916.82'
$ws.Range('E45').Value = '  +7.08%  '
$ws.Range('E46').Value = '  +5.80%  '
$ws.Range('D47').Value = '9.27'
$ws.Range('E47').Value = '  +1.17%  '
$ws.Range('D48').Value = '2.76'
$ws.Range('E48').Value = '  +2.14%  '
$ws.Range('D49').Value = '145.89'
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '3.10'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '3.10'
$ws.Range('E51').Value = '  -8.13%  '
